# Apply scheduled-runner market data updates to each sheet's Leve profit table.
# Generated from the authoritative diff: for each (sheet, row) touched, set the new
# values for columns H-N; where a column was removed by the diff (cell no longer
# present after the edit) we ClearContents() it instead of writing 0.

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 12450
$ws.Range("I19").Value = 900
$ws.Range("J19").Value = 24000
$ws.Range("K19").Value = 900
$ws.Range("L19").Value = 24000
$ws.Range("M19").Value = -725
$ws.Range("N19").Value = -24350
# Row 40
$ws.Range("H40").Value = 2231361
$ws.Range("I40").Value = 12328.5
$ws.Range("J40").Value = 6669426
$ws.Range("K40").Value = 12328.5
$ws.Range("L40").Value = 6669426
$ws.Range("M40").Value = -12153.5
$ws.Range("N40").Value = -6669776
# Row 64
$ws.Range("H64").Value = 5701
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
# Row 67
$ws.Range("H67").Value = 5701
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
# Row 74
$ws.Range("H74").Value = 35720400
$ws.Range("I74").Value = 100002840
$ws.Range("J74").Value = 7934.6665
$ws.Range("K74").Value = 100002840
$ws.Range("L74").Value = 7934.6665
$ws.Range("M74").Value = -100001904
$ws.Range("N74").Value = -9806.666499999999
# Row 77
$ws.Range("H77").Value = 35720400
$ws.Range("I77").Value = 100002840
$ws.Range("J77").Value = 7934.6665
$ws.Range("K77").Value = 500014200
$ws.Range("L77").Value = 39673.3325
$ws.Range("M77").Value = -500009520
$ws.Range("N77").Value = -49033.3325
# Row 116
$ws.Range("H116").Value = 13163257
$ws.Range("I116").Value = 50002560
$ws.Range("K116").Value = 50002560
$ws.Range("M116").Value = -49999118
# Row 127
$ws.Range("H127").Value = 4125
$ws.Range("J127").Value = 3500
$ws.Range("L127").Value = 10500
$ws.Range("N127").Value = -20420
# Row 132
$ws.Range("H132").Value = 2360.2166
$ws.Range("I132").Value = 2252.6428
$ws.Range("J132").Value = 3866.25
$ws.Range("K132").Value = 6757.928400000001
$ws.Range("L132").Value = 11598.75
$ws.Range("M132").Value = -4227.928400000001
$ws.Range("N132").Value = -16658.75

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2608211
$ws.Range("I32").Value = 2910945
$ws.Range("K32").Value = 2910945
$ws.Range("M32").Value = -2910658
# Row 74
$ws.Range("H74").Value = 61875.668
$ws.Range("J74").Value = 3839.0908
$ws.Range("L74").Value = 3839.0908
$ws.Range("N74").Value = -5587.0908
# Row 77
$ws.Range("H77").Value = 61875.668
$ws.Range("J77").Value = 3839.0908
$ws.Range("L77").Value = 19195.454
$ws.Range("N77").Value = -27931.454
# Row 97
$ws.Range("H97").Value = 4912059.5
$ws.Range("I97").Value = 791.2
$ws.Range("J97").Value = 11928158
$ws.Range("K97").Value = 791.2
$ws.Range("L97").Value = 11928158
$ws.Range("M97").Value = -295.2
$ws.Range("N97").Value = -11929150
# Row 102
$ws.Range("H102").Value = 4624.273
$ws.Range("I102").Value = 4496.4443
$ws.Range("J102").Value = 5199.5
$ws.Range("K102").Value = 4496.4443
$ws.Range("L102").Value = 5199.5
$ws.Range("M102").Value = -2874.4443
$ws.Range("N102").Value = -8443.5
# Row 110
$ws.Range("H110").Value = 19608694
$ws.Range("I110").Value = 864.75
$ws.Range("J110").Value = 66667484
$ws.Range("K110").Value = 864.75
$ws.Range("L110").Value = 66667484
$ws.Range("M110").Value = 1180.25
$ws.Range("N110").Value = -66671574
# Row 132
$ws.Range("H132").Value = 5918.8237
$ws.Range("I132").Value = 2820.7896
$ws.Range("J132").Value = 9843
$ws.Range("K132").Value = 8462.3688
$ws.Range("L132").Value = 29529
$ws.Range("M132").Value = -5932.3688
$ws.Range("N132").Value = -34589

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 48
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
# Row 86
$ws.Range("H86").Value = 47622540
$ws.Range("I86").Value = 2102.4614
$ws.Range("K86").Value = 2102.4614
$ws.Range("M86").Value = -979.4614000000001
# Row 89
$ws.Range("H89").Value = 47622540
$ws.Range("I89").Value = 2102.4614
$ws.Range("K89").Value = 10512.307
$ws.Range("M89").Value = -4896.307000000001

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 43
$ws.Range("H43").Value = 44400.25
$ws.Range("J43").Value = 44400.25
$ws.Range("L43").Value = 44400.25
$ws.Range("N43").Value = -44768.25
# Row 99
$ws.Range("H99").Value = 10702
$ws.Range("I99").Value = 16825
$ws.Range("K99").Value = 16825
$ws.Range("M99").Value = -15327
# Row 101
$ws.Range("H101").Value = 44400.25
$ws.Range("J101").Value = 44400.25
$ws.Range("L101").Value = 44400.25
$ws.Range("N101").Value = -50890.25
# Row 107
$ws.Range("H107").Value = 2102.2307
$ws.Range("I107").Value = 1565.5
$ws.Range("J107").Value = 2562.2856
$ws.Range("K107").Value = 1565.5
$ws.Range("L107").Value = 2562.2856
$ws.Range("M107").Value = 354.5
$ws.Range("N107").Value = -6402.2856
# Row 122
$ws.Range("H122").Value = 799.25
$ws.Range("I122").Value = 799
$ws.Range("K122").Value = 2397
$ws.Range("M122").Value = 53
# Row 126
$ws.Range("H126").Value = 10702
$ws.Range("I126").Value = 16825
$ws.Range("K126").Value = 50475
$ws.Range("M126").Value = -48005
# Row 132
$ws.Range("H132").Value = 4955.143
$ws.Range("I132").Value = 1925
$ws.Range("J132").Value = 10409.4
$ws.Range("K132").Value = 5775
$ws.Range("L132").Value = 31228.2
$ws.Range("M132").Value = -3245
$ws.Range("N132").Value = -36288.2

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 24
$ws.Range("H24").Value = 3560.6667
$ws.Range("I24").Value = 3454.6155
$ws.Range("K24").Value = 10363.8465
$ws.Range("M24").Value = -10133.8465
# Row 68
$ws.Range("H68").Value = 3671.25
$ws.Range("I68").Value = 1048.8
$ws.Range("J68").Value = 8042
$ws.Range("K68").Value = 3146.4
$ws.Range("L68").Value = 24126
$ws.Range("M68").Value = -2335.4
$ws.Range("N68").Value = -25748
# Row 71
$ws.Range("H71").Value = 3671.25
$ws.Range("I71").Value = 1048.8
$ws.Range("J71").Value = 8042
$ws.Range("K71").Value = 9439.199999999999
$ws.Range("L71").Value = 72378
$ws.Range("M71").Value = -5383.199999999999
$ws.Range("N71").Value = -80490
# Row 113
$ws.Range("H113").Value = 2864.3333
$ws.Range("I113").Value = 1589.5
$ws.Range("K113").Value = 4768.5
$ws.Range("M113").Value = -2598.5
# Row 121
$ws.Range("H121").Value = 16667597
$ws.Range("J121").Value = 14287115
$ws.Range("L121").Value = 42861345
$ws.Range("N121").Value = -42863965

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 3559.6667
$ws.Range("I80").Value = 3559.6667
$ws.Range("K80").Value = 3559.6667
$ws.Range("M80").Value = -2561.6667
# Row 83
$ws.Range("H83").Value = 3559.6667
$ws.Range("I83").Value = 3559.6667
$ws.Range("K83").Value = 17798.3335
$ws.Range("M83").Value = -12806.3335
# Row 97
$ws.Range("H97").Value = 1469.25
$ws.Range("I97").Value = 1245.88
$ws.Range("J97").Value = 1841.5333
$ws.Range("K97").Value = 1245.88
$ws.Range("L97").Value = 1841.5333
$ws.Range("M97").Value = -749.8800000000001
$ws.Range("N97").Value = -2833.5333
# Row 132
$ws.Range("H132").Value = 8269.4
$ws.Range("I132").Value = 3242
$ws.Range("K132").Value = 9726
$ws.Range("M132").Value = -7196

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 14292986
$ws.Range("I132").Value = 25004116
$ws.Range("K132").Value = 75012348
$ws.Range("M132").Value = -75009818

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 1007
$ws.Range("I107").Value = 760.6667
$ws.Range("K107").Value = 2282.0001
$ws.Range("M107").Value = -362.0001000000002
# Row 122
$ws.Range("H122").Value = 5338.913
$ws.Range("J122").Value = 6107
$ws.Range("L122").Value = 18321
$ws.Range("N122").Value = -23221
# Row 132
$ws.Range("H132").Value = 13908046
$ws.Range("I132").Value = 18523732
$ws.Range("J132").Value = 60990
$ws.Range("K132").Value = 55571196
$ws.Range("L132").Value = 182970
$ws.Range("M132").Value = -55568666
$ws.Range("N132").Value = -188030
